# Re-order the "Recorded By" (column G) values so that "System" (when
# present as the first, leading token of the comma-separated list) is
# moved to the end of the list instead of the beginning.
#
# Example: "System, backup@backdoor.com, system" -> "backup@backdoor.com, system, System"
#          "System, dnasr281@gmail.com"           -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -ne $val -and $val -is [string] -and $val.StartsWith("System,")) {
        $rest = $val.Substring(7).Trim()
        $cell.Value = "$rest, System"
    }
}
